$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2023" column (S) mirroring the existing year columns (B..R),
# copying the number formatting/styling from the adjacent 2022 column (R)
# and then filling in the new figures for each indicator row.

$ws.Range("R3:R14").Copy()
$ws.Range("S3:S14").PasteSpecial(-4122)

$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 1926.4
$ws.Range("S5").Value = 1929.2
$ws.Range("S6").Value = 24982
$ws.Range("S7").Value = 24520
$ws.Range("S8").Value = 1481.1
$ws.Range("S9").Value = 1068.5
$ws.Range("S10").Value = 443
$ws.Range("S11").Value = 860.8
$ws.Range("S12").Value = 240.1
$ws.Range("S13").Value = 1057.7
$ws.Range("S14").Value = 1

# Match the final selection recorded for this sheet.
$ws.Range("G21:H21").Select() | Out-Null
